# "Last edits for this term (hopefully) ^__^"
#
# The body originally held 4 paragraphs:
#   1. "Conclusion and Recommendation" (bold heading, with the _GoBack
#      bookmark sitting right before the heading run)
#   2. an empty paragraph
#   3. a long "Social networking sites ..." paragraph
#   4. a long "As for the future researchers ..." paragraph
#
# The edit collapses everything down to a single paragraph: keep the
# heading run, then move the _GoBack bookmark to just after it, followed
# by one trailing space run (23pt "body" formatting) - and drop the three
# paragraphs of body text entirely.

$d = $word.ActiveDocument

# Rebuild paragraph 1 in one shot via raw OOXML so we control the exact
# run/bookmark order (run, then bookmarkStart/bookmarkEnd, then the new
# trailing-space run) instead of the original (bookmarkStart/bookmarkEnd,
# then run).
$newHeadingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
            '<w:b/>' +
            '<w:sz w:val="36"/>' +
            '<w:szCs w:val="32"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
            '<w:b/>' +
            '<w:sz w:val="36"/>' +
            '<w:szCs w:val="32"/>' +
        '</w:rPr>' +
        '<w:t>Conclusion and Recommendation</w:t>' +
    '</w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="32"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve"> </w:t>' +
    '</w:r>' +
'</w:p>'

$d.Paragraphs(1).Range.InsertXML($newHeadingXml)

# Now remove the old paragraphs 2..N (the blank line and the two long
# body paragraphs), which all shifted down by one after the rewrite above
# but are otherwise untouched.
if ($d.Paragraphs.Count -gt 1) {
    $tailStart = $d.Paragraphs(2).Range.Start
    $tailEnd = $d.Paragraphs($d.Paragraphs.Count).Range.End
    $d.Range($tailStart, $tailEnd).Delete()
}
